$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Type" header in column H, matching the existing header style
$ws.Range("H1").Value = "Type"

# Copy the formatting from the preceding header cell (G1) onto H1 so it
# reuses the same cell style as the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Set the new column width to match the diff (~22.13 characters)
$ws.Columns.Item(8).ColumnWidth = 21.29
